$d = $word.ActiveDocument

# --- 1) Correct small transcription errors in the existing paragraph ---
$d.Content.Find.Execute("Settecento. Ma sono rimasti", $true, $false, $false, $false, $false, $true, 1, $false, "Settecento. Mi sono rimasti", 2) | Out-Null
$d.Content.Find.Execute("Eh, ma io devo contare", $true, $false, $false, $false, $false, $true, 1, $false, "Eh, ma li devo contare", 2) | Out-Null
$d.Content.Find.Execute("Che fa? Paga in contanti o con la carta? Paga Valerio. Chi? A presente Valerio, quello che ha la pescheria proprio di fronte a lei? E chi non lo conosce? Siccome deve darmi settecento euro, e allora ce li facciamo portare", $true, $false, $false, $false, $false, $true, 1, $false, "Che fa, paga in contanti o con la carta? Paga Valerio. Chi? Ha presente Valerio, quello che ha la pescheria proprio di fronte a lei. E chi non lo conosce? Siccome deve darmi settecento euro, allora ci rifacciamo portare", 2) | Out-Null
$d.Content.Find.Execute("Valerio! Eccomi. Di quei", $true, $false, $false, $false, $false, $true, 1, $false, "Valerio! Eccomi! Di quei", 2) | Out-Null
$d.Content.Find.Execute("Quattrocento sessantacinque, quattrocento sessantasei...", $true, $false, $false, $false, $false, $true, 1, $false, "Quattrocentosessantacinque, quattrocentosessantasei...", 2) | Out-Null

# --- 2) Append the Executive Summary section after the transcript paragraph ---
$d.Content.Find.Execute("Quattrocentosessantacinque, quattrocentosessantasei...", $true, $false, $false, $false, $false, $true, 1, $false, "Quattrocentosessantacinque, quattrocentosessantasei...^p^m^pExecutive Summary & Action Items^p### 📝 Riepilogo Esecutivo^p• Acquistare settecento gamberi di piccolo taglio.  ^p• Contare i gamberi con precisione, impiegando circa 15 minuti.  ^p• Acquistare una collana d’oro 18 carati a 650 euro.  ^p• Effettuare il pagamento tramite Valerio, responsabile della pescheria.  ^p• Coordinare il trasferimento di 650 euro da Valerio al gioielliere.  ^p• Garantire la consegna del denaro prima della finalizzazione dell’acquisto.  ^p§EMPTY§^p### ✅ Azioni e Responsabili^p• Valerio → Contare i settecento gamberi (data da definire)  ^p• Cliente → Acquistare collana d’oro a 650 euro (data da definire)  ^p• Valerio → Consegnare 650 euro al gioielliere Marco (data da definire)", 2) | Out-Null

# --- 3) Re-apply Heading1 style to the new section title, and Heading1 reference check ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Executive Summary & Action Items`r") {
        $p.Style = "Heading 1"
    }
}

# --- 4) Turn the placeholder line into a genuinely empty paragraph ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "§EMPTY§`r") {
        $emptyRange = $d.Range($p.Range.Start, $p.Range.End - 1)
        $emptyRange.Text = ""
    }
}
